$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF)
$data = @(
    @(7, 8),
    @(7, 9),
    @(8, 8),
    @(1, 3),
    @(8, 8),
    @(5, 6),
    @(4, 5),
    @(9, 9),
    @(7, 8),
    @(5, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
